# ---------------------------------------------------------------------------
# Applies the "Technology -> Chemistry" content rewrite plus the
# TimesNewToman -> "Times New Roman" font-name fix to the active document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Text content replacements (title, byline, email, and body copy).
# ---------------------------------------------------------------------------

Replace-Text "Technology: Bridging Humanities and Sciences" "Unveiling the Enigmatic Realm of Chemistry"
Replace-Text "Sarah Williamson" "Caroline Thompson"
Replace-Text "sarah" "caroline"
Replace-Text "williamson@eliteacademia" "thompson@gmail"
Replace-Text "org" "com"

Replace-Text "Technology, like a bridge spanning distant shores, has woven together the once disparate realms of humanities and sciences" "Chemistry, often perceived as an enigmatic ensemble of theories, equations, and experiments, unveils a captivating realm of science that seeks to unravel the very fabric of the substances that constitute our universe"
Replace-Text " This convergence has not only reshaped these fields but has also given rise to exciting new opportunities for exploration and understanding" " Chemistry is the study of matter and its properties, encompassing the transformation of substances into new substances through chemical reactions"
Replace-Text " In this essay, we delve into the dynamic interplay between technology and humanities, examining how it has transformed the study of art, literature, music, and history while simultaneously opening up innovative avenues for scientific inquiry and medical advancements" " This diverse field delves into the enigmatic mysteries of atoms, molecules, and the interactions between them, shaping the world we perceive around us"
Replace-Text "In the realm of art and literature, technology has served as a transformative muse, inspiring new forms of artistic expression and enabling the creation of immersive experiences that transcend traditional boundaries" "Chemistry permeates every aspect of our existence, from the air we breathe to the food we consume"
Replace-Text " From digital paintings that mimic the ebb and flow of brushstrokes to virtual reality installations that transport viewers to otherworldly landscapes, technology has empowered artists to push the limits of their creativity and engage audiences in unprecedented ways" " It is essential for life and plays a pivotal role in various industries, driving technological advancements and improving our quality of life"
Replace-Text " Similarly, in the realm of literature, AI-powered language models have ignited discussions about the future of narrative and opened up new possibilities for collaboration between humans and machines in the creation of literary works" " From pharmaceuticals that alleviate ailments to fertilizers that nourish crops, chemistry has become an integral part of our society, contributing immensely to healthcare, agriculture, and countless other domains"
Replace-Text "Meanwhile, in the domain of science and medicine, technology has revolutionized the way we study the body, prevent disease, and develop treatments" "The study of chemistry necessitates keen observation, analytical thinking, and a systematic approach to problem-solving"
Replace-Text " Advanced medical imaging techniques, fueled by sophisticated AI algorithms, enable us to visualize the intricate inner workings of the human body with unprecedented clarity, leading to more accurate diagnoses and targeted therapies" " It encourages curiosity, fosters creativity, and cultivates a deep appreciation for the world around us"
Replace-Text " Moreover, the development of wearable health devices and remote monitoring systems empowers individuals to take a more proactive role in managing their own health, fostering a new era of personalized and preventative medicine" " As we delve into the captivating tapestry of reactions and compounds, chemistry empowers us to understand the universe at its fundamental level, enabling us to unravel the secrets hidden within the enigmatic realms of matter"

Replace-Text "Technology has served as a catalyst for convergence between humanities and sciences, leading to profound transformations in both fields" "Chemistry, an intriguing and multifaceted realm, explores the properties and transformations of matter through chemical reactions"
Replace-Text " It has empowered artists and authors to explore new frontiers of creativity, expanded the horizons of scientific inquiry, and revolutionized medical practices" " It is a science that intricately intertwines with life, industry, and technological advancements"
Replace-Text " As technology continues to advance, we can anticipate even more remarkable breakthroughs and innovations that will redefine the very essence of humanity and science" " The study of chemistry not only enriches our understanding of the world but also cultivates critical thinking, analytical skills, and a profound appreciation for the microscopic world of atoms and molecules"

# ---------------------------------------------------------------------------
# 2. Append a new sentence (plus its own period run, matching the existing
#    run layout of one-sentence-per-run followed by a one-period-per-run)
#    to the end of the closing ("Summary") paragraph.
# ---------------------------------------------------------------------------

$summary = $d.Paragraphs.Last
$tailRange = $summary.Range
$periodRange = $d.Range($tailRange.End - 2, $tailRange.End - 1)

$sentence = " Chemistry, in essence, is an enigmatic yet captivating journey into the heart of matter and its myriad interactions"
$sentenceStart = $periodRange.End
$periodRange.InsertAfter($sentence)
$sentenceRange = $d.Range($sentenceStart, $sentenceStart + $sentence.Length)
$sentenceRange.Font.Name = "Times New Roman"
$sentenceRange.Font.Color = 0

$fullStop = "."
$fullStopStart = $sentenceRange.End
$sentenceRange.InsertAfter($fullStop)
$fullStopRange = $d.Range($fullStopStart, $fullStopStart + $fullStop.Length)
$fullStopRange.Font.Name = "Times New Roman"
$fullStopRange.Font.Color = 0

# ---------------------------------------------------------------------------
# 3. Add a new, empty trailing paragraph after the closing paragraph.
# ---------------------------------------------------------------------------

$d.Paragraphs.Last.Range.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 4. Fix the font name everywhere: "TimesNewToman" -> "Times New Roman".
#    Iterate paragraph-by-paragraph (excluding each paragraph mark) so the
#    rename only touches existing run properties instead of also stamping a
#    new <w:pPr><w:rPr> on every paragraph.
# ---------------------------------------------------------------------------

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    if ($r.Start -lt $r.End) {
        $r.Font.Name = "Times New Roman"
    }
}

Write-Output "edit complete"
